$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for price cells whose refreshed values would
# otherwise be auto-converted to numbers by Excel (e.g. "1.00" -> 1).
$textForceRows = @(4,5,6,8,9,10,11,12,13,14,18,20,21,22,23,24,25,27,28,29,30,31,32,33,34,36,37,38,40,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textForceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Apply cell content updates per the refreshed crypto data snapshot
# Row 2
$ws.Range('D2').Value = '67.328.28'
$ws.Range('E2').Value = '  -1.60%  '

# Row 3
$ws.Range('D3').Value = '3.600.23'
$ws.Range('E3').Value = '  -2.62%  '

# Row 4
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '578.90'
$ws.Range('E5').Value = '  -5.34%  '

# Row 6
$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '193.18'
$ws.Range('E6').Value = '  -0.91%  '

# Row 7
$ws.Range('D7').Value = '3.596.78'
$ws.Range('E7').Value = '  -2.55%  '

# Row 8
$ws.Range('D8').Value = '0.620'
$ws.Range('E8').Value = '  -2.29%  '

# Row 9
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  -0.13%  '

# Row 10
$ws.Range('D10').Value = '0.682'
$ws.Range('E10').Value = '  -5.74%  '

# Row 11
$ws.Range('D11').Value = '0.152'
$ws.Range('E11').Value = '  -5.02%  '

# Row 12
$ws.Range('D12').Value = '56.04'
$ws.Range('E12').Value = '  -5.97%  '

# Row 13
$ws.Range('D13').Value = '0.0000278'
$ws.Range('E13').Value = '  -3.07%  '

# Row 14
$ws.Range('D14').Value = '9.91'
$ws.Range('E14').Value = '  -4.99%  '

# Row 15
$ws.Range('D15').Value = '4.169.00'
$ws.Range('E15').Value = '  -2.94%  '

# Row 16
$ws.Range('D16').Value = '3.596.58'
$ws.Range('E16').Value = '  -2.97%  '

# Row 17
$ws.Range('E17').Value = '  -1.34%  '

# Row 18
$ws.Range('D18').Value = '18.46'
$ws.Range('E18').Value = '  -4.88%  '

# Row 19
$ws.Range('D19').Value = '67.215.05'
$ws.Range('E19').Value = '  -1.63%  '

# Row 20
$ws.Range('D20').Value = '12.23'
$ws.Range('E20').Value = '  -4.59%  '

# Row 21
$ws.Range('D21').Value = '1.07'
$ws.Range('E21').Value = '  -6.68%  '

# Row 22
$ws.Range('D22').Value = '402.73'
$ws.Range('E22').Value = '  -1.35%  '

# Row 23
$ws.Range('D23').Value = '4.23'
$ws.Range('E23').Value = '  -8.53%  '

# Row 24
$ws.Range('D24').Value = '86.12'
$ws.Range('E24').Value = '  -4.20%  '

# Row 25
$ws.Range('D25').Value = '11.40'
$ws.Range('E25').Value = '  -0.99%  '

# Row 26
$ws.Range('E26').Value = '  -3.77%  '

# Row 27
$ws.Range('D27').Value = '12.53'
$ws.Range('E27').Value = '  -3.95%  '

# Row 28
$ws.Range('D28').Value = '6.10'
$ws.Range('E28').Value = '  +1.06%  '

# Row 29
$ws.Range('D29').Value = '3.66'
$ws.Range('E29').Value = '  -2.70%  '

# Row 30
$ws.Range('D30').Value = '9.00'
$ws.Range('E30').Value = '  -6.18%  '

# Row 31
$ws.Range('D31').Value = '7.66'
$ws.Range('E31').Value = '  -1.47%  '

# Row 32
$ws.Range('D32').Value = '31.32'
$ws.Range('E32').Value = '  -4.19%  '

# Row 33
$ws.Range('D33').Value = '634.20'
$ws.Range('E33').Value = '  +0.07%  '

# Row 34
$ws.Range('D34').Value = '12.23'
$ws.Range('E34').Value = '  -3.67%  '

# Row 35
$ws.Range('E35').Value = '  -5.16%  '

# Row 36
$ws.Range('D36').Value = '64.13'
$ws.Range('E36').Value = '  -4.97%  '

# Row 37
$ws.Range('D37').Value = '42.56'
$ws.Range('E37').Value = '  -11.35%  '

# Row 38
$ws.Range('D38').Value = '0.402'
$ws.Range('E38').Value = '  -2.68%  '

# Row 39
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0798'
$ws.Range('E39').Value = '  -2.35%  '

# Row 40
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.04%  '

# Row 41
$ws.Range('D41').Value = '3.181.01'
$ws.Range('E41').Value = '  +8.49%  '

# Row 42
$ws.Range('D42').Value = '0.134'
$ws.Range('E42').Value = '  -3.20%  '

# Row 43
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '2.74'
$ws.Range('E43').Value = '  +4.78%  '

# Row 44
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  -0.28%  '

# Row 45
$ws.Range('D45').Value = '2.99'
$ws.Range('E45').Value = '  -1.51%  '

# Row 46
$ws.Range('D46').Value = '0.0418'
$ws.Range('E46').Value = '  -5.73%  '

# Row 47
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.66'
$ws.Range('E47').Value = '  -0.88%  '

# Row 48
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.131'
$ws.Range('E48').Value = '  -6.30%  '

# Row 49
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '3.08'
$ws.Range('E49').Value = '  +1.34%  '

# Row 50
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '141.62'
$ws.Range('E50').Value = '  -3.05%  '

# Row 51
$ws.Range('D51').Value = '8.60'
$ws.Range('E51').Value = '  -7.77%  '
